$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 504, shifting existing rows 504:537 down to 506:539
$ws.Rows.Item(504).Resize(2).Insert()

# New row 504: Pimiento / Zafiro rojo / Región de Arica y Parinacota
$ws.Range("A504").Value = 5
$ws.Range("B504").Value = "Macroferia Regional de Talca"
$ws.Range("C504").Value = "Maule"
$ws.Range("D504").Value = 44714
$ws.Range("E504").Value = 7
$ws.Range("F504").Value = 100112002
$ws.Range("G504").Value = "Pimiento"
$ws.Range("H504").Value = "Zafiro rojo"
$ws.Range("I504").Value = "Primera"
$ws.Range("J504").Value = 200
$ws.Range("K504").Value = 43000
$ws.Range("L504").Value = 43000
$ws.Range("M504").Value = 43000
$ws.Range("N504").Value = "$/caja 15 kilos"
$ws.Range("O504").Value = "Región de Arica y Parinacota"
$ws.Range("P504").Value = 2867
$ws.Range("Q504").Value = 15
$ws.Range("R504").Value = "Hortaliza"

# New row 505: Pimiento / Zafiro verde / Región de Arica y Parinacota
$ws.Range("A505").Value = 5
$ws.Range("B505").Value = "Macroferia Regional de Talca"
$ws.Range("C505").Value = "Maule"
$ws.Range("D505").Value = 44714
$ws.Range("E505").Value = 7
$ws.Range("F505").Value = 100112002
$ws.Range("G505").Value = "Pimiento"
$ws.Range("H505").Value = "Zafiro verde"
$ws.Range("I505").Value = "Primera"
$ws.Range("J505").Value = 200
$ws.Range("K505").Value = 25000
$ws.Range("L505").Value = 25000
$ws.Range("M505").Value = 25000
$ws.Range("N505").Value = "$/caja 15 kilos"
$ws.Range("O505").Value = "Región de Arica y Parinacota"
$ws.Range("P505").Value = 1667
$ws.Range("Q505").Value = 15
$ws.Range("R505").Value = "Hortaliza"
